$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.577.27"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "3.102.47"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.16%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.102.36"
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("E12").Value = "  +2.40%  "
$ws.Range("D13").Value = "3.633.62"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "57.634.09"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "3.100.99"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "336.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +2.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("E28").Value = "  +2.39%  "
$ws.Range("E29").Value = "  +1.51%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("E32").Value = "  +1.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "156.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("E36").Value = "  +2.69%  "
$ws.Range("E37").Value = "  +2.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("E39").Value = "  +2.71%  "
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.94"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("B42").Value = "RenzoRestakedETH"
$ws.Range("C42").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D42").Value = "3.139.98"
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("E43").Value = "  +4.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Value = "2.302.65"
$ws.Range("E47").Value = "  +1.23%  "
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.980"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.39%  "
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("E51").Value = "  +2.19%  "
